$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Rename sheet from ParentSoils to the generic Actions concept
$ws.Name = "Actions"

# Add the new OpportunityCost column header
$ws.Cells.Item(1, 6).Value = "OpportunityCost"
$ws.Cells.Item(1, 6).Style = "Normal"

# Populate the OpportunityCost values for rows 2-17
$values = @(
    0,
    11969.116740092626,
    6205.1723024022067,
    317.83533591054601,
    72373.744154991917,
    2236.700270436133,
    16368.774399517197,
    1224.2486701957214,
    0,
    3240.2366959041756,
    0,
    3030.2103241448754,
    0,
    2124.2628833762665,
    2.3919495589002002,
    2925.5337950916137
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $values[$i]
    $cell.Style = "Normal"
}

# Widen the new column
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666

# Update the active selection on the sheet
$ws.Range("H8").Select() | Out-Null
